# Apply cryptos list update (GitHub Actions style scrape refresh)

function Set-TextValue($Range, $Text) {
    # Force the cell to store $Text as literal text (matches original inlineStr cells),
    # instead of letting Excel auto-convert numeric-looking strings to numbers,
    # while preserving the original (default) cell style.
    $savedStyle = $Range.Style
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.Style = $savedStyle
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
Set-TextValue $ws.Range("D2") '23.130.19'
$ws.Range("E2").Value = '  -3.10%  '

# Row 3
Set-TextValue $ws.Range("D3") '1.607.37'
$ws.Range("E3").Value = '  -2.44%  '

# Row 4
Set-TextValue $ws.Range("D4") '1.002'
$ws.Range("E4").Value = '  -0.10%  '

# Row 5
Set-TextValue $ws.Range("D5") '1.001'
$ws.Range("E5").Value = '  -0.11%  '

# Row 6
Set-TextValue $ws.Range("D6") '302.22'
$ws.Range("E6").Value = '  -2.16%  '

# Row 7
Set-TextValue $ws.Range("D7") '0.3768'
$ws.Range("E7").Value = '  -3.06%  '

# Row 8
Set-TextValue $ws.Range("D8") '0.3641'
$ws.Range("E8").Value = '  -4.70%  '

# Row 9
Set-TextValue $ws.Range("D9") '48.76'
$ws.Range("E9").Value = '  -4.81%  '

# Row 10
Set-TextValue $ws.Range("D10") '1.002'
$ws.Range("E10").Value = '  -0.09%  '

# Row 11
$ws.Range("E11").Value = '  -6.14%  '

# Row 12
Set-TextValue $ws.Range("D12") '0.08056'
$ws.Range("E12").Value = '  -4.24%  '

# Row 13
Set-TextValue $ws.Range("D13") '22.91'
$ws.Range("E13").Value = '  -4.01%  '

# Row 15
Set-TextValue $ws.Range("D15") '7.605'
$ws.Range("E15").Value = '  -3.39%  '

# Row 16
Set-TextValue $ws.Range("D16") '0.00001259'
$ws.Range("E16").Value = '  -4.17%  '

# Row 17
Set-TextValue $ws.Range("D17") '1.602.71'
$ws.Range("E17").Value = '  -2.77%  '

# Row 18
Set-TextValue $ws.Range("D18") '91.39'
$ws.Range("E18").Value = '  -3.10%  '

# Row 19
Set-TextValue $ws.Range("D19") '0.06783'
$ws.Range("E19").Value = '  -2.88%  '

# Row 20
Set-TextValue $ws.Range("D20") '18.29'
$ws.Range("E20").Value = '  -6.78%  '

# Row 21
Set-TextValue $ws.Range("D21") '6.556'
$ws.Range("E21").Value = '  -5.21%  '

# Row 22
Set-TextValue $ws.Range("D22") '1.001'
$ws.Range("E22").Value = '  -0.13%  '

# Row 23
Set-TextValue $ws.Range("D23") '13.09'
$ws.Range("E23").Value = '  -4.19%  '

# Row 24
Set-TextValue $ws.Range("D24") '23.158.14'
$ws.Range("E24").Value = '  -3.03%  '

# Row 25
Set-TextValue $ws.Range("D25") '2.357'
$ws.Range("E25").Value = '  -3.78%  '

# Row 26
Set-TextValue $ws.Range("D26") '2.860'
$ws.Range("E26").Value = '  -3.81%  '

# Row 27
Set-TextValue $ws.Range("D27") '21.05'
$ws.Range("E27").Value = '  -4.56%  '

# Row 28
Set-TextValue $ws.Range("D28") '150.32'
$ws.Range("E28").Value = '  -0.24%  '

# Row 29
Set-TextValue $ws.Range("D29") '5.249'
$ws.Range("E29").Value = '  -2.58%  '

# Row 30
Set-TextValue $ws.Range("D30") '132.22'
$ws.Range("E30").Value = '  -4.61%  '

# Row 31
Set-TextValue $ws.Range("D31") '2.390'
$ws.Range("E31").Value = '  -4.41%  '

# Row 32
Set-TextValue $ws.Range("D32") '6.803'
$ws.Range("E32").Value = '  -12.58%  '

# Row 33
Set-TextValue $ws.Range("D33") '1.778.76'
$ws.Range("E33").Value = '  -2.77%  '

# Row 34
Set-TextValue $ws.Range("D34") '0.9691'
$ws.Range("E34").Value = '  -7.31%  '

# Row 35
Set-TextValue $ws.Range("D35") '0.07699'
$ws.Range("E35").Value = '  -4.12%  '

# Row 36
Set-TextValue $ws.Range("D36") '0.02774'
$ws.Range("E36").Value = '  -5.93%  '

# Row 37
$ws.Range("B37").Value = 'Algorand'
$ws.Range("C37").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws.Range("D37") '0.2543'
$ws.Range("E37").Value = '  -4.98%  '

# Row 38
$ws.Range("B38").Value = 'InternetComputer(DFINITY)'
$ws.Range("C38").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue $ws.Range("D38") '6.211'
$ws.Range("E38").Value = '  -7.33%  '

# Row 39
Set-TextValue $ws.Range("D39") '10.12'
$ws.Range("E39").Value = '  -6.40%  '

# Row 40
Set-TextValue $ws.Range("D40") '0.08840'
$ws.Range("E40").Value = '  -2.76%  '

# Row 41
Set-TextValue $ws.Range("D41") '1.388'
$ws.Range("E41").Value = '  -2.28%  '

# Row 42
Set-TextValue $ws.Range("D42") '0.7153'
$ws.Range("E42").Value = '  -5.24%  '

# Row 43
Set-TextValue $ws.Range("D43") '12.78'
$ws.Range("E43").Value = '  -4.69%  '

# Row 44
Set-TextValue $ws.Range("D44") '15.79'
$ws.Range("E44").Value = '  -3.08%  '

# Row 45
Set-TextValue $ws.Range("D45") '0.6600'
$ws.Range("E45").Value = '  -4.77%  '

# Row 46
Set-TextValue $ws.Range("D46") '1.000'
$ws.Range("E46").Value = '  -0.14%  '

# Row 47
Set-TextValue $ws.Range("D47") '2.291'
$ws.Range("E47").Value = '  -6.58%  '

# Row 48
Set-TextValue $ws.Range("D48") '3.974'
$ws.Range("E48").Value = '  -2.50%  '

# Row 49
Set-TextValue $ws.Range("D49") '0.07988'
$ws.Range("E49").Value = '  -3.43%  '

# Row 50
Set-TextValue $ws.Range("D50") '131.61'
$ws.Range("E50").Value = '  -1.61%  '

# Row 51
Set-TextValue $ws.Range("D51") '1.169'
$ws.Range("E51").Value = '  -3.01%  '

